$d = $word.ActiveDocument

# 1. "In prokaryotes, there is only transcription." -> add "(i.e. no splicing)"
$d.Content.Find.Execute(
    "there is only transcription.", $true, $false, $false, $false, $false,
    $true, 1, $false, "there is only transcription (i.e. no splicing).", 2) | Out-Null

# 2. "Intron - Noncoding portion of DNA between exons." -> "Non-coding"
$d.Content.Find.Execute(
    "Noncoding portion of DNA between exons.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Non-coding portion of DNA between exons.", 2) | Out-Null

# 3. "Goal of Bioinformatics - Enable the discover of new..." -> "discovery"
$d.Content.Find.Execute(
    "Enable the discover of new", $true, $false, $false, $false, $false,
    $true, 1, $false, "Enable the discovery of new", 2) | Out-Null

# 4. "Pathogen: Disease causing microbe." -> "Disease causing agent (e.g. microbe)."
$d.Content.Find.Execute(
    " Disease causing microbe.", $true, $false, $false, $false, $false,
    $true, 1, $false, " Disease causing agent (e.g. microbe).", 2) | Out-Null

# 5. "PAM (Percent Accepted Mutation)" -> "PAM (Point Accepted Mutation)"
$d.Content.Find.Execute(
    "PAM (Percent Accepted Mutation)", $true, $false, $false, $false, $false,
    $true, 1, $false, "PAM (Point Accepted Mutation)", 2) | Out-Null

# 6. "Accepted Mutation: Any mutation that is not fatal to the organism or destroy the protein."
#    -> "...organism nor does it destroy the protein."
$d.Content.Find.Execute(
    "Any mutation that is not fatal to the organism or destroy the protein.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Any mutation that is not fatal to the organism nor does it destroy the protein.", 2) | Out-Null

# 7. "...homologs and not due to change." -> "...and not due to chance."
$d.Content.Find.Execute(
    " and not due to change.", $true, $false, $false, $false, $false,
    $true, 1, $false, " and not due to chance.", 2) | Out-Null

# 8. Fill in the two empty table cells (the "Spliceosome / 30% Identity Score / (empty) / (empty)" row)
$t = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Range.Text -like "*Spliceosome*30% Identity Score*") {
        $t = $candidate
    }
}

$cell3 = $t.Cell(1, 3)
$cell3.Range.Text = "UCSC Genome Browser – Tool used to visualize genome information."
$cell3Start = $t.Cell(1, 3).Range.Start
$termLen3 = "UCSC Genome Browser".Length
$termRange3 = $d.Range($cell3Start, $cell3Start + $termLen3)
$termRange3.Font.Bold = 1
$termRange3.Font.Color = 16711680

$cell4 = $t.Cell(1, 4)
$cell4.Range.Text = "Complement – Indication that the coding sequence is on the complementary strand in GenBank."
$cell4Start = $t.Cell(1, 4).Range.Start
$termLen4 = "Complement ".Length
$termRange4 = $d.Range($cell4Start, $cell4Start + $termLen4)
$termRange4.Font.Bold = 1
$termRange4.Font.Color = 16711680

# 9. Footer page number field cached result: 9 -> 5
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
foreach ($f in $ftr.Range.Fields) {
    $f.Result.Text = "5"
}
